$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing style of the Price column, force it to Text so that
# numeric-looking strings (e.g. "1.001") are written back as literal text
# instead of being auto-converted into numbers by Excel, then restore style.
$priceRange = $ws.Range("D2:D51")
$origPriceStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.645.81'
$ws.Range("E2").Value = '  -2.22%  '

$ws.Range("D3").Value = '1.846.97'
$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '314.65'
$ws.Range("E5").Value = '  -1.44%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.17%  '

$ws.Range("D7").Value = '0.4257'
$ws.Range("E7").Value = '  -2.80%  '

$ws.Range("D8").Value = '0.3659'
$ws.Range("E8").Value = '  -1.60%  '

$ws.Range("D9").Value = '45.44'
$ws.Range("E9").Value = '  +0.57%  '

$ws.Range("D10").Value = '0.07272'
$ws.Range("E10").Value = '  -3.44%  '

$ws.Range("D11").Value = '0.9012'
$ws.Range("E11").Value = '  -4.02%  '

$ws.Range("D12").Value = '20.75'
$ws.Range("E12").Value = '  -2.64%  '

$ws.Range("D13").Value = '1.828.96'
$ws.Range("E13").Value = '  -2.12%  '

$ws.Range("D14").Value = '5.369'
$ws.Range("E14").Value = '  -1.94%  '

$ws.Range("D15").Value = '6.574'
$ws.Range("E15").Value = '  -2.36%  '

$ws.Range("D16").Value = '0.06850'
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").Value = '78.08'
$ws.Range("E18").Value = '  -4.98%  '

$ws.Range("D19").Value = '0.000008811'
$ws.Range("E19").Value = '  -3.24%  '

$ws.Range("E20").Value = '  -0.18%  '

$ws.Range("D21").Value = '15.47'
$ws.Range("E21").Value = '  -3.39%  '

$ws.Range("D22").Value = '27.626.12'
$ws.Range("E22").Value = '  -2.26%  '

$ws.Range("D23").Value = '4.969'
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("D24").Value = '10.62'
$ws.Range("E24").Value = '  -1.24%  '

$ws.Range("D25").Value = '2.061.55'
$ws.Range("E25").Value = '  -1.38%  '

$ws.Range("D26").Value = '2.045'
$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("D27").Value = '154.02'
$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("D28").Value = '18.29'
$ws.Range("E28").Value = '  -0.53%  '

$ws.Range("D29").Value = '5.275'
$ws.Range("E29").Value = '  -1.58%  '

$ws.Range("D30").Value = '1.830'
$ws.Range("E30").Value = '  +5.54%  '

$ws.Range("D31").Value = '110.77'
$ws.Range("E31").Value = '  -3.08%  '

$ws.Range("D32").Value = '0.08896'
$ws.Range("E32").Value = '  -1.83%  '

$ws.Range("D33").Value = '0.7716'
$ws.Range("E33").Value = '  -4.01%  '

$ws.Range("D34").Value = '4.557'
$ws.Range("E34").Value = '  -6.16%  '

$ws.Range("D35").Value = '2.968'
$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("D36").Value = '1.087'
$ws.Range("E36").Value = '  -7.35%  '

$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("D38").Value = '0.05412'
$ws.Range("E38").Value = '  -1.02%  '

$ws.Range("E39").Value = '  -2.60%  '

$ws.Range("D40").Value = '0.01929'
$ws.Range("E40").Value = '  -1.14%  '

$ws.Range("D41").Value = '2.920'
$ws.Range("E41").Value = '  -2.02%  '

$ws.Range("D42").Value = '0.5069'
$ws.Range("E42").Value = '  -3.55%  '

$ws.Range("D43").Value = '6.820'
$ws.Range("E43").Value = '  -4.62%  '

$ws.Range("E44").Value = '  -1.94%  '

$ws.Range("D45").Value = '8.248'
$ws.Range("E45").Value = '  -6.20%  '

$ws.Range("D46").Value = '0.06634'
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.4724'
$ws.Range("E47").Value = '  -3.22%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '10.33'
$ws.Range("E48").Value = '  -2.42%  '

$ws.Range("D49").Value = '105.57'
$ws.Range("E49").Value = '  -2.10%  '

$ws.Range("D50").Value = '1.000'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("D51").Value = '1.639'
$ws.Range("E51").Value = '  -2.62%  '

# Restore the original style/format for the Price column
$priceRange.Style = $origPriceStyle
